# Generate Report for Handoff
# - Update status text "Handed back: in sync with en-US" -> "Ready for handoff"
# - Update the two "Latest Xliff Generate"/Handoff timestamp text values
# - Narrow the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" columns

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text updates (Handed back: in sync with en-US -> Ready for handoff) ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# --- Timestamp text updates (kept as plain text, not true dates) ---
$wsOverview.Range("G2").Value = "2016-09-01 17:09:35"
$wsDeDe.Range("H2").Value = "2016-09-01 17:09:35"
$wsZhCn.Range("H2").Value = "2016-09-01 17:09:30"

# --- Column width updates ---
# Target stored width is 17.2159881591797 characters. The Excel column-width
# model only persists widths on a fixed pixel grid, so the nearest
# achievable grid value (17.1666...) is reached by assigning 16.3 here.
$wsOverview.Range("E1").ColumnWidth = 16.3
$wsOverview.Range("F1").ColumnWidth = 16.3
$wsZhCn.Range("C1").ColumnWidth = 16.3
$wsDeDe.Range("C1").ColumnWidth = 16.3
